# "book list 2024 update" -- append newly-read books to the "Library" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library")

# New rows to append (Book Name, Author, Date, Category), matching the
# existing table layout (columns A-D, date column styled like the rest of
# column C).
$newRows = @(
    @{ Row = 76; Book = "Essentialism";                                         Author = "Greg McKeown";        Date = "2024-09-05"; Category = "Selfhelp" },
    @{ Row = 77; Book = "The Good Success";                                      Author = "Dr Emmanuel Mango";   Date = "2024-06-04"; Category = "Selfhelp" },
    @{ Row = 78; Book = "Failing Forward";                                       Author = "John C Maxwell";      Date = "2024-07-21"; Category = "Selfhelp" },
    @{ Row = 79; Book = "The Prize: The Epic Quest for Oil, Money & Power";      Author = "Daniel Yergin";       Date = "2023-12-12"; Category = "Business" },
    @{ Row = 80; Book = "Why (Not) Me: Memoir";                                  Author = "John C Gichinga";     Date = "2024-01-12"; Category = "Biography" },
    @{ Row = 81; Book = "One Night at the Call Centre";                          Author = "Chetan Bhagat";       Date = "2022-01-12"; Category = "Fiction" },
    @{ Row = 82; Book = "The flame trees of Thika";                             Author = " Elspeth Huxley";     Date = "2022-01-12"; Category = "Fiction" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Book
    $ws.Range("B$row").Value = $r.Author
    $ws.Range("C$row").Value = [double]([datetime]::ParseExact($r.Date, "yyyy-MM-dd", $null).ToOADate())
    $ws.Range("C$row").NumberFormat = "[$-409]d\-mmm\-yyyy;@"
    $ws.Range("D$row").Value = $r.Category
}

# Reflect the new selection / scroll position from the authored workbook.
$ws.Activate()
$ws.Range("C82").Select() | Out-Null

# Cosmetic theme rename the author made while editing (best-effort -- not
# all hosts persist theme metadata, but harmless to attempt).
try {
    $theme = $wb.Theme
    $theme.ThemeColorScheme.Name = "Office 2013 - 2022"
    $theme.ThemeFontScheme.Name = "Office 2013 - 2022"
    $theme.Name = "Office 2013 - 2022 Theme"
} catch {
}
